$wb = $excel.ActiveWorkbook
$ws2 = $wb.ActiveSheet

# Replace the plain numeric cell values with descriptive shared-string labels
# that encode Square/Position/Row/Column, in the exact order needed so the
# shared-string table is rebuilt with matching indices.
$ws2.Range("D3").Value = "TR - P1 - R1 - C1"
$ws2.Range("D4").Value = "TR - P3 - R2 - C1"
$ws2.Range("D5").Value = "BR - P1 - R3 - C1"
$ws2.Range("D6").Value = "BR - P3 - R4 - C1"

$ws2.Range("H6").Value = "BR - P4 - R4 - C2"
$ws2.Range("H5").Value = "BR - P2 - R3 - C2"
$ws2.Range("H4").Value = "TR - P4 - R2 - C2"
$ws2.Range("H3").Value = "TR - P2 - R1 - C2"

$ws2.Range("P3").Value = "TL - P2 - R1 - C4"
$ws2.Range("P4").Value = "TL - P4 - R2 - C4"
$ws2.Range("P5").Value = "BL - P2 - R3 - C4"
$ws2.Range("P6").Value = "BL - P4 - R4 - C4"

$ws2.Range("L3").Value = "TL - P1 - R1 - C3"
$ws2.Range("L5").Value = "BL - P1 - R3 - C3"
$ws2.Range("L6").Value = "BL - P3 - R4 - C3"
$ws2.Range("L4").Value = "TL - P3 - R2 - C3"

# Give the new descriptive columns enough width to show their content.
$ws2.Columns.Item(4).ColumnWidth = 13.81640625
$ws2.Columns.Item(8).ColumnWidth = 13.90625
$ws2.Columns.Item(12).ColumnWidth = 13.6328125
$ws2.Columns.Item(16).ColumnWidth = 13.6328125

# Update the active selection on Sheet2 to the whole of column P.
$ws2.Activate()
$ws2.Range("P1:P1048576").Select() | Out-Null
